$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <= original row 4 data
$ws.Cells.Item(2, 4).Value = 44545
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 550
$ws.Cells.Item(2, 11).Value = 1700
$ws.Cells.Item(2, 12).Value = 1800
$ws.Cells.Item(2, 13).Value = 1755
$ws.Cells.Item(2, 14).Value = "$/kilo"
$ws.Cells.Item(2, 15).Value = "Provincia de Linares"
$ws.Cells.Item(2, 16).Value = 1755

# Row 3 <= original row 10 data
$ws.Cells.Item(3, 4).Value = 44511
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 600
$ws.Cells.Item(3, 11).Value = 1300
$ws.Cells.Item(3, 12).Value = 1400
$ws.Cells.Item(3, 13).Value = 1350
$ws.Cells.Item(3, 14).Value = "$/kilo"
$ws.Cells.Item(3, 15).Value = "Provincia de Linares"
$ws.Cells.Item(3, 16).Value = 1350

# Row 4 <= original row 5 data
$ws.Cells.Item(4, 4).Value = 44510
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 600
$ws.Cells.Item(4, 11).Value = 1300
$ws.Cells.Item(4, 12).Value = 1400
$ws.Cells.Item(4, 13).Value = 1350
$ws.Cells.Item(4, 14).Value = "$/kilo"
$ws.Cells.Item(4, 15).Value = "Provincia de Linares"
$ws.Cells.Item(4, 16).Value = 1350

# Row 5 <= original row 7 data
$ws.Cells.Item(5, 4).Value = 44526
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 1500
$ws.Cells.Item(5, 12).Value = 1600
$ws.Cells.Item(5, 13).Value = 1550
$ws.Cells.Item(5, 14).Value = "$/kilo"
$ws.Cells.Item(5, 15).Value = "Provincia de Linares"
$ws.Cells.Item(5, 16).Value = 1550

# Row 6 <= original row 12 data
$ws.Cells.Item(6, 4).Value = 44489
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 600
$ws.Cells.Item(6, 11).Value = 1400
$ws.Cells.Item(6, 12).Value = 1500
$ws.Cells.Item(6, 13).Value = 1450
$ws.Cells.Item(6, 14).Value = "$/kilo"
$ws.Cells.Item(6, 15).Value = "Provincia de Linares"
$ws.Cells.Item(6, 16).Value = 1450

# Row 7 <= original row 6 data
$ws.Cells.Item(7, 4).Value = 44876
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 350
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 1600
$ws.Cells.Item(7, 13).Value = 1557
$ws.Cells.Item(7, 14).Value = "$/kilo"
$ws.Cells.Item(7, 15).Value = "Provincia de Linares"
$ws.Cells.Item(7, 16).Value = 1557

# Row 8 <= original row 18 data
$ws.Cells.Item(8, 4).Value = 44860
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 1100
$ws.Cells.Item(8, 11).Value = 1500
$ws.Cells.Item(8, 12).Value = 1700
$ws.Cells.Item(8, 13).Value = 1609
$ws.Cells.Item(8, 14).Value = "$/kilo"
$ws.Cells.Item(8, 15).Value = "Provincia de Linares"
$ws.Cells.Item(8, 16).Value = 1609

# Row 9 <= original row 11 data
$ws.Cells.Item(9, 4).Value = 44477
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 500
$ws.Cells.Item(9, 11).Value = 1400
$ws.Cells.Item(9, 12).Value = 1500
$ws.Cells.Item(9, 13).Value = 1460
$ws.Cells.Item(9, 14).Value = "$/kilo"
$ws.Cells.Item(9, 15).Value = "Provincia de Linares"
$ws.Cells.Item(9, 16).Value = 1460

# Row 10 <= original row 19 data
$ws.Cells.Item(10, 4).Value = 44496
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 550
$ws.Cells.Item(10, 11).Value = 1500
$ws.Cells.Item(10, 12).Value = 2000
$ws.Cells.Item(10, 13).Value = 1773
$ws.Cells.Item(10, 14).Value = "$/paquete"
$ws.Cells.Item(10, 15).Value = "Provincia de Linares"
$ws.Cells.Item(10, 16).Value = 1773

# Row 11 <= original row 2 data
$ws.Cells.Item(11, 4).Value = 44875
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 300
$ws.Cells.Item(11, 11).Value = 1500
$ws.Cells.Item(11, 12).Value = 1600
$ws.Cells.Item(11, 13).Value = 1550
$ws.Cells.Item(11, 14).Value = "$/kilo"
$ws.Cells.Item(11, 15).Value = "Provincia de Linares"
$ws.Cells.Item(11, 16).Value = 1550

# Row 12 <= original row 3 data
$ws.Cells.Item(12, 4).Value = 44839
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 500
$ws.Cells.Item(12, 11).Value = 1700
$ws.Cells.Item(12, 12).Value = 1800
$ws.Cells.Item(12, 13).Value = 1760
$ws.Cells.Item(12, 14).Value = "$/kilo"
$ws.Cells.Item(12, 15).Value = "Provincia de Linares"
$ws.Cells.Item(12, 16).Value = 1760

# Row 13 <= original row 14 data
$ws.Cells.Item(13, 4).Value = 44868
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 1000
$ws.Cells.Item(13, 11).Value = 1200
$ws.Cells.Item(13, 12).Value = 1300
$ws.Cells.Item(13, 13).Value = 1250
$ws.Cells.Item(13, 14).Value = "$/kilo"
$ws.Cells.Item(13, 15).Value = "Región del Maule"
$ws.Cells.Item(13, 16).Value = 1250

# Row 14 <= original row 15 data
$ws.Cells.Item(14, 4).Value = 44868
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Segunda"
$ws.Cells.Item(14, 10).Value = 200
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 12).Value = 1000
$ws.Cells.Item(14, 13).Value = 1000
$ws.Cells.Item(14, 14).Value = "$/kilo"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 1000

# Row 15 <= original row 8 data
$ws.Cells.Item(15, 4).Value = 44524
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 1500
$ws.Cells.Item(15, 12).Value = 1600
$ws.Cells.Item(15, 13).Value = 1550
$ws.Cells.Item(15, 14).Value = "$/kilo"
$ws.Cells.Item(15, 15).Value = "Provincia de Talca"
$ws.Cells.Item(15, 16).Value = 1550

# Row 18 <= original row 13 data
$ws.Cells.Item(18, 4).Value = 44468
$ws.Cells.Item(18, 8).Value = "Verde"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 1800
$ws.Cells.Item(18, 12).Value = 2000
$ws.Cells.Item(18, 13).Value = 1920
$ws.Cells.Item(18, 14).Value = "$/kilo"
$ws.Cells.Item(18, 15).Value = "Provincia de Linares"
$ws.Cells.Item(18, 16).Value = 1920

# Row 19 <= original row 9 data
$ws.Cells.Item(19, 4).Value = 44519
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 250
$ws.Cells.Item(19, 11).Value = 1200
$ws.Cells.Item(19, 12).Value = 1300
$ws.Cells.Item(19, 13).Value = 1240
$ws.Cells.Item(19, 14).Value = "$/kilo"
$ws.Cells.Item(19, 15).Value = "Provincia de Linares"
$ws.Cells.Item(19, 16).Value = 1240

